$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new row (row 6) describing the "COMP_BRANCH" instruction encoding ---
# The row mirrors the layout of the existing GPR_NUM row (row 4) / header row
# (row 2): opcode | Rs1 | Rs2 | 0 | 0 | 0 | Address

# 1) Merge the destination groups first - doing this before formatting/values
#    keeps the existing cell-style indices intact instead of Excel re-splitting
#    the border formatting across the merged range.
$ws.Range("B6:F6").Merge() | Out-Null
$ws.Range("G6:J6").Merge() | Out-Null
$ws.Range("K6:N6").Merge() | Out-Null
$ws.Range("R6:AG6").Merge() | Out-Null

# 2) Copy the cell formatting (one source cell at a time) onto the new row so
#    that every destination cell reuses the matching existing style.
# B6:F6 "opcode" field -> same formatting as B4:F4
for ($c = 2; $c -le 6; $c++) {
    $ws.Cells.Item(4, $c).Copy() | Out-Null
    $ws.Cells.Item(6, $c).PasteSpecial(-4122)
}

# G6:J6 "Rs1" field -> same formatting as K2:N2
$srcCols = @(11, 12, 13, 14)
$dstCols = @(7, 8, 9, 10)
for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $ws.Cells.Item(2, $srcCols[$i]).Copy() | Out-Null
    $ws.Cells.Item(6, $dstCols[$i]).PasteSpecial(-4122)
}

# K6:N6 "Rs2" field -> same formatting as O2:R2
$srcCols = @(15, 16, 17, 18)
$dstCols = @(11, 12, 13, 14)
for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $ws.Cells.Item(2, $srcCols[$i]).Copy() | Out-Null
    $ws.Cells.Item(6, $dstCols[$i]).PasteSpecial(-4122)
}

# R6:AF6 "Address" field -> same formatting as P4:AD4; AG6 -> same as AG4
$srcCols = 16..30
$dstCols = 18..32
for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $ws.Cells.Item(4, $srcCols[$i]).Copy() | Out-Null
    $ws.Cells.Item(6, $dstCols[$i]).PasteSpecial(-4122)
}
$ws.Cells.Item(4, 33).Copy() | Out-Null
$ws.Cells.Item(6, 33).PasteSpecial(-4122)

# 3) Fill in the cell values. "Address" is written before "COMP_BRANCH" so the
#    shared-string table grows in the same order as the saved workbook.
$ws.Cells.Item(6, 2).Value = "opcode"
$ws.Cells.Item(6, 7).Value = "Rs1"
$ws.Cells.Item(6, 11).Value = "Rs2"
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 0
$ws.Cells.Item(6, 18).Value = "Address"
$ws.Cells.Item(6, 1).Value = "COMP_BRANCH"

# --- Update the sheet selection to match the post-edit state ---
$ws.Range("A6").Select() | Out-Null
